$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" column (G) values replacing former "Strike#" values, rows 2-16
$newValues = @(0, 1, 0, 0, 0, 1, 1, 1, 0, 2, 1, 1, 1, 2, 2)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $newValues[$i]
}
